# Log Week 17 data into the OFF and DEF target depth sheets.

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 222
$wsOff.Range("C2").Value = 165
$wsOff.Range("D2").Value = 60
$wsOff.Range("E2").Value = 25
$wsOff.Range("G2").Value = 5

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 220
$wsDef.Range("C2").Value = 148
$wsDef.Range("D2").Value = 70
$wsDef.Range("E2").Value = 34
$wsDef.Range("F2").Value = 3
$wsDef.Range("G2").Value = 2
